$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row index: one past the current last row (row 4 -> row 5)
$newRow = 5

# Plain text values (non numeric-looking) can be written directly; Excel
# keeps them as text cells the same way the existing rows are stored.
$ws.Cells.Item($newRow, 2).Value = "أحمد شريم"
$ws.Cells.Item($newRow, 4).Value = "الصمود"
$ws.Cells.Item($newRow, 5).Value = "الرحلة 3"
$ws.Cells.Item($newRow, 6).Value = "C1"
$ws.Cells.Item($newRow, 7).Value = "UNICEF"
$ws.Cells.Item($newRow, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:١٥:٣٢ م"

# A5 ("23") and C5 ("234") look like numbers, so a plain .Value assignment
# would be auto-coerced to a Number cell. Enter them as a string-formula
# result on a scratch cell, then copy/paste-special the *value* into place
# so the destination keeps its text typing without picking up a new
# NumberFormat/style (which a direct NumberFormat="@" + Value assignment
# would otherwise stamp onto the cell).
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""23"""
$scratch.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)

$scratch.Formula = "=""234"""
$scratch.Copy()
$ws.Cells.Item($newRow, 3).PasteSpecial(-4163)

$scratch.Clear()
